$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$lot2052 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"

$ws.Range("B24").Value = $lot2028
$ws.Range("C24").Value = $lot2028

$ws.Range("B25").Value = $lot2052
$ws.Range("C25").Value = $lot2052
